$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 206.88889
$ws.Range("I33").Value = 234.85715
$ws.Range("J33").Value = 109
$ws.Range("K33").Value = 234.85715
$ws.Range("L33").Value = 109
$ws.Range("M33").Value = -5.85714999999999
$ws.Range("N33").Value = -567
$ws.Range("H70").Value = 2794.7693
$ws.Range("I70").Value = 1495.25
$ws.Range("K70").Value = 4485.75
$ws.Range("M70").Value = -4215.75
$ws.Range("H73").Value = 2794.7693
$ws.Range("I73").Value = 1495.25
$ws.Range("K73").Value = 4485.75
$ws.Range("M73").Value = -3549.75
$ws.Range("H86").Value = 90913460
$ws.Range("I86").Value = 333337150
$ws.Range("K86").Value = 333337150
$ws.Range("M86").Value = -333336027
$ws.Range("H89").Value = 90913460
$ws.Range("I89").Value = 333337150
$ws.Range("K89").Value = 1666685750
$ws.Range("M89").Value = -1666680134
$ws.Range("H100").Value = 1812
$ws.Range("I100").Value = 1812
$ws.Range("K100").Value = 1812
$ws.Range("M100").Value = -1271
$ws.Range("H112").Value = 1547.75
$ws.Range("I112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("M112").ClearContents()
$ws.Range("H131").Value = 597193.5
$ws.Range("I131").Value = 778368.25
$ws.Range("K131").Value = 2335104.75
$ws.Range("M131").Value = -2330064.75
$ws.Range("H138").Value = 2824.6272
$ws.Range("I138").Value = 1184
$ws.Range("J138").Value = 3010.3584
$ws.Range("K138").Value = 3552
$ws.Range("L138").Value = 9031.075199999999
$ws.Range("M138").Value = 1588
$ws.Range("N138").Value = -19311.0752

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 143.6
$ws.Range("I5").Value = 157
$ws.Range("K5").Value = 157
$ws.Range("M5").Value = -45
$ws.Range("H32").Value = 2404.8262
$ws.Range("I32").Value = 2404.8262
$ws.Range("K32").Value = 2404.8262
$ws.Range("M32").Value = -2117.8262

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 143.6
$ws.Range("I4").Value = 157
$ws.Range("K4").Value = 157
$ws.Range("M4").Value = -42
$ws.Range("H15").Value = 300
$ws.Range("I15").Value = 300
$ws.Range("K15").Value = 300
$ws.Range("M15").Value = -73
$ws.Range("H22").Value = 803.2727
$ws.Range("I22").Value = 729.75
$ws.Range("J22").Value = 999.3333
$ws.Range("K22").Value = 729.75
$ws.Range("L22").Value = 999.3333
$ws.Range("M22").Value = -556.75
$ws.Range("N22").Value = -1345.3333
$ws.Range("H134").Value = 2408.8438
$ws.Range("I134").Value = 1666.1052
$ws.Range("K134").Value = 4998.3156
$ws.Range("M134").Value = -2463.3156

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 3464
$ws.Range("I7").Value = 4758.636
$ws.Range("K7").Value = 4758.636
$ws.Range("M7").Value = -4645.636
$ws.Range("H16").Value = 1999
$ws.Range("I16").Value = 1999
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1999
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -1712
$ws.Range("H99").Value = 3179.4
$ws.Range("I99").Value = 2474.25
$ws.Range("J99").Value = 6000
$ws.Range("K99").Value = 2474.25
$ws.Range("L99").Value = 6000
$ws.Range("M99").Value = -976.25
$ws.Range("N99").Value = -8996
$ws.Range("H113").Value = 1999
$ws.Range("I113").Value = 1999
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1999
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = 171
$ws.Range("H126").Value = 3179.4
$ws.Range("I126").Value = 2474.25
$ws.Range("J126").Value = 6000
$ws.Range("K126").Value = 7422.75
$ws.Range("L126").Value = 18000
$ws.Range("M126").Value = -4952.75
$ws.Range("N126").Value = -22940
$ws.Range("H132").Value = 12826140
$ws.Range("I132").Value = 3568.9333
$ws.Range("K132").Value = 10706.7999
$ws.Range("M132").Value = -8176.7999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 72.616165
$ws.Range("I4").Value = 72.616165
$ws.Range("K4").Value = 217.848495
$ws.Range("M4").Value = -105.848495
$ws.Range("H33").Value = 2887.4285
$ws.Range("I33").Value = 41
$ws.Range("J33").Value = 6682.6665
$ws.Range("K33").Value = 246
$ws.Range("L33").Value = 40095.999
$ws.Range("M33").Value = 37
$ws.Range("N33").Value = -40661.999
$ws.Range("H87").Value = 4210
$ws.Range("I87").Value = 1525
$ws.Range("K87").Value = 4575
$ws.Range("M87").Value = -3327
$ws.Range("H90").Value = 4210
$ws.Range("I90").Value = 1525
$ws.Range("K90").Value = 13725
$ws.Range("M90").Value = -7485
$ws.Range("H125").Value = 5000
$ws.Range("J125").Value = 5000
$ws.Range("L125").Value = 15000
$ws.Range("N125").Value = -24840
$ws.Range("H131").Value = 11906907
$ws.Range("I131").Value = 27779854
$ws.Range("J131").Value = 2196.4167
$ws.Range("K131").Value = 83339562
$ws.Range("L131").Value = 6589.250100000001
$ws.Range("M131").Value = -83334522
$ws.Range("N131").Value = -16669.2501

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 5000
$ws.Range("J33").Value = 5000
$ws.Range("L33").Value = 5000
$ws.Range("N33").Value = -5504
$ws.Range("H102").Value = 1999.5
$ws.Range("I102").Value = 1995
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 1995
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = -373
$ws.Range("N102").Value = -5244
$ws.Range("H132").Value = 2645.7896
$ws.Range("I132").Value = 2519.2856
$ws.Range("K132").Value = 7557.8568
$ws.Range("M132").Value = -5027.8568

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2260.625
$ws.Range("I7").Value = 2239
$ws.Range("K7").Value = 2239
$ws.Range("M7").Value = -2127
$ws.Range("H35").Value = 2566
$ws.Range("I35").Value = 2888
$ws.Range("J35").Value = 1600
$ws.Range("K35").Value = 2888
$ws.Range("L35").Value = 1600
$ws.Range("M35").Value = -2552
$ws.Range("N35").Value = -2272
$ws.Range("H40").Value = 1250
$ws.Range("I40").Value = 1250
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 1250
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -1114
$ws.Range("H68").Value = 5000
$ws.Range("I68").Value = 5000
$ws.Range("K68").Value = 5000
$ws.Range("M68").Value = -4251
$ws.Range("H71").Value = 5000
$ws.Range("I71").Value = 5000
$ws.Range("K71").Value = 25000
$ws.Range("M71").Value = -21256
$ws.Range("H93").Value = 1140.6666
$ws.Range("I93").Value = 465.66666
$ws.Range("K93").Value = 465.66666
$ws.Range("M93").Value = 782.33334
$ws.Range("H100").Value = 4983
$ws.Range("I100").Value = 4499.75
$ws.Range("K100").Value = 4499.75
$ws.Range("M100").Value = -3958.75
$ws.Range("H126").Value = 2260.625
$ws.Range("I126").Value = 2239
$ws.Range("K126").Value = 6717
$ws.Range("M126").Value = -4247
$ws.Range("H132").Value = 10611.305
$ws.Range("I132").Value = 6553.1113
$ws.Range("J132").Value = 13220.143
$ws.Range("K132").Value = 19659.3339
$ws.Range("L132").Value = 39660.429
$ws.Range("M132").Value = -17129.3339
$ws.Range("N132").Value = -44720.429
$ws.Range("H136").Value = 6353.95
$ws.Range("I136").Value = 5233.3335
$ws.Range("J136").Value = 7270.8184
$ws.Range("K136").Value = 15700.0005
$ws.Range("L136").Value = 21812.4552
$ws.Range("M136").Value = -13150.0005
$ws.Range("N136").Value = -26912.4552

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6333.3335
$ws.Range("J62").Value = 5000
$ws.Range("L62").Value = 5000
$ws.Range("N62").Value = -6248
$ws.Range("H65").Value = 6333.3335
$ws.Range("J65").Value = 5000
$ws.Range("L65").Value = 25000
$ws.Range("N65").Value = -31240
$ws.Range("H113").Value = 1613.5
$ws.Range("I113").Value = 1613.5
$ws.Range("K113").Value = 4840.5
$ws.Range("M113").Value = -2670.5
$ws.Range("H126").Value = 2642.5557
$ws.Range("I126").Value = 2683.8572
$ws.Range("K126").Value = 8051.571599999999
$ws.Range("M126").Value = -5581.571599999999
$ws.Range("H132").Value = 2035.7097
$ws.Range("I132").Value = 1925.7778
$ws.Range("J132").Value = 2777.75
$ws.Range("K132").Value = 5777.3334
$ws.Range("L132").Value = 8333.25
$ws.Range("M132").Value = -3247.3334
$ws.Range("N132").Value = -13393.25
$ws.Range("H136").Value = 427906.12
$ws.Range("I136").Value = 2673.087
$ws.Range("J136").Value = 1825100.4
$ws.Range("K136").Value = 8019.261
$ws.Range("L136").Value = 5475301.199999999
$ws.Range("M136").Value = -5469.261
$ws.Range("N136").Value = -5480401.199999999
